$wb = $excel.ActiveWorkbook

# --- Rename the second sheet ---
$ws = $wb.Worksheets.Item("testSubscribe")
$ws.Name = "testWorkFrontJobs"

# --- Resize the recorded window geometry (cosmetic bookView; best effort) ---
try { $excel.ActiveWindow.Width = 27630 } catch {}
try { $excel.ActiveWindow.Height = 11700 } catch {}

# --- Rebuild the header/data for the WorkFront Jobs fields (column by column) ---
$ws.Range("C1").Value = "Client Name"
$ws.Range("C2").Value = "Automation test client"

$ws.Range("D1").Value = "NAN"
$ws.Range("D2").Value = 25693698

$ws.Range("E1").Value = "EIN"
$ws.Range("E2").Value = 256987458

$ws.Range("F1").Value = "EIN Type"
$ws.Range("F2").Value = "Tax Filing - Federal"

$ws.Range("G1").Value = "Address"
$ws.Range("H1").Value = "City"
$ws.Range("I1").Value = "State"
$ws.Range("J1").Value = "Zip"

$ws.Range("G2").Value = "Street 119"
$ws.Range("H2").Value = "Dallas"
$ws.Range("I2").Value = "TX"
$ws.Range("J2").Value = 75898

# --- Column width tweaks to fit the new header text (best-fit-like autosize) ---
$ws.Columns.Item(3).ColumnWidth = 20.333333333333332
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 17

# --- View/selection tweaks ---
$ws.Range("C4").Select()
